$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates per the commit diff.
# Numeric-looking text values use NumberFormat="@" + Style="Normal" to keep them
# stored as text (matching the source data, which keeps prices as text strings)
# without leaving a residual custom style on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '245.77'
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '24.10'
$ws.Range("D3").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.272'
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05779'
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.466'
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.126'
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8167'
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8535'
$ws.Range("D9").Style = "Normal"

$ws.Range("B10").Value = 'One'

$ws.Range("C10").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.009855'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = '9OneONEBestin24h'

$ws.Range("B11").Value = 'WazirX'

$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1359'
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = '10WazirXWRX'

$ws.Range("B12").Value = 'MandalaExchangeToken'

$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.06923'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = '11MandalaExchangeTokenMDX'

$ws.Range("B13").Value = 'LiechtensteinCryptoassetsExchange'

$ws.Range("C13").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03125'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '12LiechtensteinCryptoassetsExchangeLCX'

$ws.Range("B14").Value = 'BitrueCoin'

$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.02916'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = '13BitrueCoinBTR'

$ws.Range("B15").Value = 'BitMartToken'

$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.09379'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = '14BitMartTokenBMX'

$ws.Range("B16").Value = 'MCDex'

$ws.Range("C16").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.735'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '15MCDexMCB'

$ws.Range("B17").Value = 'BitForexToken'

$ws.Range("C17").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.001515'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '16BitForexTokenBF'

$ws.Range("B18").Value = 'CoinExToken'

$ws.Range("C18").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.04666'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = '17CoinExTokenCET'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.006199'
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.001235'
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.004624'
$ws.Range("D21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.00006886'
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.503'
$ws.Range("D23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.148'
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.3194'
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1320'
$ws.Range("D26").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1359'
$ws.Range("D27").Style = "Normal"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0002327'
$ws.Range("D28").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03651'
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006217'
$ws.Range("D41").Style = "Normal"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.003393'
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008387'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = '43LocalTradersLCT'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005246'
$ws.Range("D45").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000748'
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.3692'
$ws.Range("D47").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002273'
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002096'
$ws.Range("D49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0001996'
$ws.Range("D50").Style = "Normal"
